$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the contact-info runs into a single run and add the extra
#    "github.com/danielmartincraig" + " • linkedin.com/danielcraig23" text.
#    Matching across the existing run boundaries (and the spellStart/spellEnd
#    proofErr wrapper) collapses them into one clean run, exactly like the
#    target markup.
# ---------------------------------------------------------------------------
$bullet = [char]0x2022
$oldContact = "(803)389-6750 " + $bullet + " danielmartincraig@gmail.com " + $bullet + " github.com/danielmartincraig " + $bullet + " linkedin.com/danielcraig23"
$newContact = "(803)389-6750 " + $bullet + " danielmartincraig@gmail.com " + $bullet + " github.com/danielmartincraig " + $bullet + " linkedin.com/danielcraig23"
$d.Content.Find.Execute($oldContact, $true, $false, $false, $false, $false, $true, 1, $false, $newContact, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Merge the "Web Engineering I and II" runs (removing the gramStart/
#    gramEnd proofErr wrapper around "Engineering") into a single run.
# ---------------------------------------------------------------------------
$webEng = "•    Web Engineering I and II"
$d.Content.Find.Execute($webEng, $true, $false, $false, $false, $false, $true, 1, $false, $webEng, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark from the end of the document (after
#    "Fluent in Spanish") into the new OBJECTIVE paragraph we are about to
#    add. Delete the old one now; it gets re-added once the new paragraph's
#    text exists.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4. Insert the new "OBJECTIVE:" heading paragraph right after the contact
#    info paragraph.
# ---------------------------------------------------------------------------
$paraIndex = 0
$contactParaIndex = -1
foreach ($p in $d.Paragraphs) {
    $paraIndex = $paraIndex + 1
    if ($p.Range.Text -like "*linkedin.com/danielcraig23*") {
        $contactParaIndex = $paraIndex
    }
}

$contactPara = $d.Paragraphs.Item($contactParaIndex)
$insertionPoint = $d.Range($contactPara.Range.End, $contactPara.Range.End)
$insertionPoint.InsertAfter([char]13)

$objPara = $d.Paragraphs.Item($contactParaIndex + 1)
$objPara.Style = "Heading1"

$runObjective = $d.Range($objPara.Range.Start, $objPara.Range.Start)
$runObjective.InsertAfter("OBJECTIVE: ")

$runEager = $d.Range($runObjective.End, $runObjective.End)
$runEager.InsertAfter("Eager to drive back-end solutions at ")
$runEager.Font.Size = 12

$runRincon = $d.Range($runEager.End, $runEager.End)
$runRincon.InsertAfter("Rincon")
$runRincon.Font.Size = 12

$bookmarkPos = $runRincon.End

$runBasis = $d.Range($runRincon.End, $runRincon.End)
$runBasis.InsertAfter(" on a full-time basis")
$runBasis.Font.Size = 12

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
